# Edit: Add ML inference Lambda with optimized layer (218MB unzipped)

$wb = $excel.ActiveWorkbook

# --- Sheet references (tab order matches xl/_rels/workbook.xml.rels) ---
$wsOverview = $wb.Worksheets.Item(1)   # Project Overview
$wsPhase4   = $wb.Worksheets.Item(5)   # Phase 4 - ML Integration
$wsChangeLog = $wb.Worksheets.Item(9)  # Change Log

# =========================================================================
# 1. Change Log: record Phase 2 completion and Phase 3 completion entries
#    (duplicate the most recent row's formatting/content, then customize).
# =========================================================================
$wsChangeLog.Range("A6:F6").Copy()
$wsChangeLog.Range("A7:F8").PasteSpecial(-4104)
$excel.CutCopyMode = $false

$wsChangeLog.Range("A7").Value = 45987
$wsChangeLog.Range("A8").Value = 45998

$wsChangeLog.Range("C8").Value = "Phase 3"

$wsChangeLog.Range("E7").Value = "API & compute layer complete - 100%"
$wsChangeLog.Range("E8").Value = "Multi-Modal Activity Tracking complete - 97%"

$wsChangeLog.Range("F8").Value = "Medium"

$wsChangeLog.Range("A7").Borders.Item(9).LineStyle = 1
$wsChangeLog.Range("A7").Borders.Item(9).Weight = 2
$wsChangeLog.Range("B7:F7").Borders.Item(9).LineStyle = 1
$wsChangeLog.Range("B7:F7").Borders.Item(9).Weight = 2

$wsChangeLog.Range("A6").Borders.Item(9).LineStyle = -4142
$wsChangeLog.Range("B6:F6").Borders.Item(9).LineStyle = -4142

$wsChangeLog.Range("D10").Select()

# =========================================================================
# 2. Phase 4 - ML Integration: mark the "Create Lambda Handler for ML
#    Inference" task group (rows 9-12) complete, log actual hours, and
#    record the blocker + mitigation note hit while packaging it.
# =========================================================================
$wsPhase4.Range("E9").Value = "Complete"
$wsPhase4.Range("H9").Value = 2

$wsPhase4.Range("E10").Value = "Complete"
$wsPhase4.Range("H10").Value = 1.5

$wsPhase4.Range("E11").Value = "Complete"
$wsPhase4.Range("H11").Value = 0.5

$wsPhase4.Range("E12").Value = "Complete"
$wsPhase4.Range("H12").Value = 4
$wsPhase4.Range("I12").Value = "AWS Resource size limit"
$wsPhase4.Range("J12").Value = "Lambda Layers to separate dependencies"

# =========================================================================
# 3. Project Overview: refresh KEY METRICS now that four more tasks are
#    complete, and leave a note on the Phase 4 summary row.
# =========================================================================
$wsOverview.Range("B25").Value = "70 (51%)"
$wsOverview.Range("B28").Value = "63 (46%)"
$wsOverview.Range("B31").Value = "76.5 / 120"

$wsOverview.Range("H18").Value = "Lambda size limit - optimization"

$wsOverview.Select()
$wsOverview.Range("A8:H8").Select()
$excel.ActiveWindow.ScrollRow = 8
$wsOverview.Range("I20").Select()
